$wb = $excel.ActiveWorkbook

# Map of Sheet name -> Cell -> New Value for the "Elapsed Duration(Hrs)" column (G)
$updates = @{
    "R1" = @{
        "G2" = "3951:43:18"
        "G3" = "91:15:56"
        "G4" = "114:15:56"
    }
    "R2" = @{
        "G2" = "12133:06:59"
        "G3" = "3262:50:28"
        "G4" = "501:02:02"
    }
    "R4" = @{
        "G2" = "2978:56:48"
        "G3" = "206:09:03"
        "G4" = "94:21:28"
        "G5" = "91:59:01"
    }
    "R5" = @{
        "G2" = "452:55:47"
    }
    "R6" = @{
        "G2" = "93:28:05"
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $updates[$sheetName]
    foreach ($cellRef in $cells.Keys) {
        $ws.Range($cellRef).Value = $cells[$cellRef]
    }
}
